# Pipe network analysis: replace the 6x6 diameter table with a smaller
# 4x4 table of newly measured diameters.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Drop the rows/columns that no longer belong to the (now smaller) table.
$ws.Range("F1:G1").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E3:E4").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("A6:G7").ClearContents()

# Update the surviving diameter values for the new table.
$ws.Range("B3").Value = 0.9
$ws.Range("B4").Value = 0.9
$ws.Range("C4").Value = 0.9
$ws.Range("C5").Value = 0.9
$ws.Range("D5").Value = 0.9

# Match the author's final selection on the new table's corner cell.
$ws.Range("D5").Select()
